# Update "想去人数" (want-to-go count) values in column F
# for the sheets "展览" and "全部类型" — these two sheets carry identical
# content, and both need the same refreshed snapshot numbers.

$wb = $excel.ActiveWorkbook

# row -> new F value
$updates = @{
    2  = 1768
    5  = 1133
    6  = 42
    7  = 12214
    10 = 487
    13 = 887
    14 = 13576
    15 = 13690
    17 = 160
    19 = 39
    20 = 1014
    21 = 101
    23 = 2149
    24 = 199
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
